# Updated cryptos list on Mon Sep  2 18:51:48 UTC 2024 with GitHub Actions
# Applies per-row Price (D) and Volume(1h) (E) updates to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "58.149.65"; E = "  -0.32%  " },
    @{ Row = 3; D = "2.498.06"; E = "  +0.81%  " },
    @{ Row = 4; D = "1.00"; E = "  +0.34%  " },
    @{ Row = 5; D = "520.33"; E = "  -0.28%  " },
    @{ Row = 6; D = "132.69"; E = "  -1.45%  " },
    @{ Row = 7; D = "1.00"; E = "  +0.72%  " },
    @{ Row = 8; D = "0.560"; E = "  +0.18%  " },
    @{ Row = 9; D = "2.515.72"; E = "  +1.09%  " },
    @{ Row = 10; D = "0.0976"; E = "  -0.92%  " },
    @{ Row = 11; D = $null; E = "  -1.54%  " },
    @{ Row = 12; D = $null; E = "  -3.54%  " },
    @{ Row = 13; D = "0.330"; E = "  -2.79%  " },
    @{ Row = 14; D = "2.949.62"; E = "  +1.20%  " },
    @{ Row = 15; D = "58.184.80"; E = "  -0.11%  " },
    @{ Row = 16; D = "22.05"; E = "  -0.59%  " },
    @{ Row = 17; D = "0.0000134"; E = "  -0.76%  " },
    @{ Row = 18; D = "2.509.80"; E = "  +1.33%  " },
    @{ Row = 19; D = "10.63"; E = "  -0.59%  " },
    @{ Row = 20; D = "321.16"; E = "  +0.02%  " },
    @{ Row = 21; D = "4.14"; E = "  -1.30%  " },
    @{ Row = 22; D = "6.15"; E = "  +7.20%  " },
    @{ Row = 23; D = "0.998"; E = "  -0.18%  " },
    @{ Row = 24; D = "64.58"; E = "  +0.09%  " },
    @{ Row = 25; D = "0.406"; E = "  -1.29%  " },
    @{ Row = 26; D = $null; E = "  +0.28%  " },
    @{ Row = 27; D = $null; E = "  -0.72%  " },
    @{ Row = 28; D = "7.37"; E = "  -0.55%  " },
    @{ Row = 29; D = "0.0₃0750"; E = "  -0.33%  " },
    @{ Row = 30; D = $null; E = "  +1.01%  " },
    @{ Row = 31; D = "167.57"; E = "  -1.02%  " },
    @{ Row = 32; D = $null; E = "  -0.08%  " },
    @{ Row = 33; D = "6.30"; E = "  -0.19%  " },
    @{ Row = 34; D = "0.998"; E = "  -0.02%  " },
    @{ Row = 35; D = "0.994"; E = "  -0.13%  " },
    @{ Row = 36; D = "18.14"; E = "  +0.20%  " },
    @{ Row = 37; D = $null; E = "  -6.03%  " },
    @{ Row = 38; D = "3.90"; E = "  -2.94%  " },
    @{ Row = 39; D = $null; E = "  +0.46%  " },
    @{ Row = 40; D = "36.36"; E = "  -0.77%  " },
    @{ Row = 41; D = "0.768"; E = "  -4.09%  " },
    @{ Row = 42; D = "275.50"; E = "  +0.00%  " },
    @{ Row = 43; D = "3.43"; E = "  -0.83%  " },
    @{ Row = 44; D = "5.00"; E = "  -3.72%  " },
    @{ Row = 45; D = "129.20"; E = "  +3.84%  " },
    @{ Row = 46; D = "0.597"; E = "  -0.15%  " },
    @{ Row = 47; D = "0.0920"; E = "  +0.96%  " },
    @{ Row = 48; D = "0.0500"; E = "  +1.76%  " },
    @{ Row = 49; D = "17.67"; E = "  -0.61%  " },
    @{ Row = 50; D = "0.0213"; E = "  -0.32%  " },
    @{ Row = 51; D = "16.85"; E = "  -1.25%  " }

)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Force text storage so numeric-looking strings (e.g. "1.00",
        # "520.33") keep their exact original formatting instead of being
        # auto-converted to numbers by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
